$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C67").Value = "US Core Medication Profile | US Core MedicationRequest Profile"
$ws.Range("E67").Value = "US Core Medication Profile|  US Core MedicationRequest Profile"

$ws.Range("C68").Value = "US Core MedicationRequest Profile"
$ws.Range("E68").Value = "US Core MedicationRequest Profile"

$ws.Range("C69").Value = "US Core MedicationRequest Profile"
$ws.Range("E69").Value = "US Core MedicationRequest Profile"

$ws.Range("C70").Value = "US Core MedicationRequest Profile"
$ws.Range("E70").Value = "US Core MedicationRequest Profile"

$ws.Range("C71").Value = "US Core MedicationDispense Profile"
$ws.Range("E71").Value = "US Core Medication Profile|  US Core MedicationDispense Profile"
